$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'297.22"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-1.19%"
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'31.38"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-0.20%"
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'5.106"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.98%"
$ws.Range("E4").Style = "Normal"

$ws.Range("E5").Value = "'8.93%"
$ws.Range("E5").Style = "Normal"

$ws.Range("E6").Value = "'37.97%"
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'7.766"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-0.84%"
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'3.794"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'1.22%"
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.9304"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-0.03%"
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.1773"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'4.26%"
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.07306"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'3.92%"
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.09004"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'11.20%"
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'0.03022"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.36%"
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'0.1004"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.92%"
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.001500"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.62%"
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'0.005862"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-4.75%"
$ws.Range("E16").Style = "Normal"

$ws.Range("E17").Value = "'1.92%"
$ws.Range("E17").Style = "Normal"

$ws.Range("E18").Value = "'1.19%"
$ws.Range("E18").Style = "Normal"

$ws.Range("E19").Value = "'-0.27%"
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'0.1347"
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").Value = "'3.380"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-25.87%"
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'0.1651"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'4.33%"
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'0.04589"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-1.15%"
$ws.Range("E23").Style = "Normal"

$ws.Range("E24").Value = "'2.35%"
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'0.004411"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-7.10%"
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'0.0001197"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-7.77%"
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'0.0003429"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'83.20%"
$ws.Range("E27").Style = "Normal"

$ws.Range("D39").Value = "'0.01756"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'2.17%"
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'0.04470"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-0.98%"
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.006873"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-3.16%"
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'0.1342"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'0.07%"
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'0.002140"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-2.54%"
$ws.Range("E43").Style = "Normal"

$ws.Range("E44").Value = "'-8.49%"
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.00006648"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'6.58%"
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'0.03%"
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'0.008754"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-14.29%"
$ws.Range("E47").Style = "Normal"

$ws.Range("E48").Value = "'-57.27%"
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'0.00002101"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'0.03%"
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'0.0002001"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.10%"
$ws.Range("E50").Style = "Normal"
